$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value would otherwise be auto-detected as a number by Excel;
# force them to remain Text (matching the workbook's original inline-string / text cells).
$textFormatCells = @("D8", "D10", "D14", "D18", "D20", "D22", "D25", "D26", "D27", "D31", "D32", "D37", "D39", "D40", "D43", "D48", "D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.919.06"
$ws.Range("E2").Value = "  +1.03%  "

$ws.Range("D3").Value = "1.640.27"
$ws.Range("E3").Value = "  +0.32%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("E5").Value = "  +0.56%  "

$ws.Range("E6").Value = "  +0.39%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "23.65"
$ws.Range("E8").Value = "  +1.28%  "

$ws.Range("E9").Value = "  -0.91%  "

$ws.Range("D10").Value = "0.0614"
$ws.Range("E10").Value = "  +0.35%  "

$ws.Range("E11").Value = "  +0.51%  "

$ws.Range("D12").Value = "1.871.41"
$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("D13").Value = "1.641.28"
$ws.Range("E13").Value = "  +0.34%  "

$ws.Range("D14").Value = "4.09"
$ws.Range("E14").Value = "  +1.04%  "

$ws.Range("E15").Value = "  +3.55%  "

$ws.Range("D17").Value = "27.909.77"
$ws.Range("E17").Value = "  +1.07%  "

$ws.Range("D18").Value = "230.90"
$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("E19").Value = "  +0.41%  "

$ws.Range("D20").Value = "7.60"
$ws.Range("E20").Value = "  +0.11%  "

$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D22").Value = "10.76"
$ws.Range("E22").Value = "  +1.02%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("E24").Value = "  -3.41%  "

$ws.Range("D25").Value = "151.50"
$ws.Range("E25").Value = "  +1.19%  "

$ws.Range("D26").Value = "6.94"
$ws.Range("E26").Value = "  +0.52%  "

$ws.Range("D27").Value = "15.68"
$ws.Range("E27").Value = "  +0.79%  "

$ws.Range("E28").Value = "  +0.14%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  +0.94%  "

$ws.Range("D31").Value = "0.0484"
$ws.Range("E31").Value = "  -0.23%  "

$ws.Range("D32").Value = "3.33"
$ws.Range("E32").Value = "  +1.42%  "

$ws.Range("E33").Value = "  +0.42%  "

$ws.Range("D34").Value = "1.392.50"
$ws.Range("E34").Value = "  -5.92%  "

$ws.Range("E35").Value = "  +1.31%  "

$ws.Range("E36").Value = "  +0.64%  "

$ws.Range("D37").Value = "0.889"
$ws.Range("E37").Value = "  +0.63%  "

$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").Value = "0.919"
$ws.Range("E39").Value = "  -2.42%  "

$ws.Range("D40").Value = "0.555"
$ws.Range("E40").Value = "  -0.81%  "

$ws.Range("E41").Value = "  -0.61%  "

$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("D43").Value = "66.20"
$ws.Range("E43").Value = "  -3.00%  "

$ws.Range("E44").Value = "  +4.30%  "

$ws.Range("E45").Value = "  +1.33%  "

$ws.Range("E46").Value = "  -0.33%  "

$ws.Range("D47").Value = "1.780.80"
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("D48").Value = "87.96"
$ws.Range("E48").Value = "  +0.27%  "

$ws.Range("E49").Value = "  +0.69%  "

$ws.Range("E50").Value = "  +0.38%  "

$ws.Range("D51").Value = "7.60"
$ws.Range("E51").Value = "  -1.59%  "
